$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.092547051789356
$ws.Range("D2").Value = 1.095061633883406
$ws.Range("E2").Value = 1.094370488614894
$ws.Range("F2").Value = 1.105850107209807
$ws.Range("I2").Value = 1.074865683611512
$ws.Range("J2").Value = 1.09736728733586
$ws.Range("K2").Value = 1.09769739405178
$ws.Range("L2").Value = 1.097008004927287
$ws.Range("M2").Value = 1.108458761183432
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.093903339520339
$ws.Range("D3").Value = 1.09618666833401
$ws.Range("E3").Value = 1.095583711910469
$ws.Range("F3").Value = 1.10709306257443
$ws.Range("I3").Value = 1.075402019186793
$ws.Range("J3").Value = 1.098387544671605
$ws.Range("K3").Value = 1.098642669306878
$ws.Range("L3").Value = 1.098041134009831
$ws.Range("M3").Value = 1.109523644434394
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.094779930220305
$ws.Range("D4").Value = 1.096913706579135
$ws.Range("E4").Value = 1.096367522697387
$ws.Range("F4").Value = 1.107896475731042
$ws.Range("I4").Value = 1.075747285865876
$ws.Range("J4").Value = 1.099046155134506
$ws.Range("K4").Value = 1.099252796655644
$ws.Range("L4").Value = 1.098707835193246
$ws.Range("M4").Value = 1.110211258595817
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.095148209904217
$ws.Range("D5").Value = 1.097219133103953
$ws.Range("E5").Value = 1.096696747282982
$ws.Range("F5").Value = 1.108234028007825
$ws.Range("I5").Value = 1.075892012335355
$ws.Range("J5").Value = 1.099322664036074
$ws.Range("K5").Value = 1.09950893134995
$ws.Range("L5").Value = 1.09898768861921
$ws.Range("M5").Value = 1.110499991235378
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.095210031802221
$ws.Range("D6").Value = 1.097270402735818
$ws.Range("E6").Value = 1.09675200868072
$ws.Range("F6").Value = 1.108290692748272
$ws.Range("I6").Value = 1.075916287764307
$ws.Range("J6").Value = 1.099369069448314
$ws.Range("K6").Value = 1.099551916307367
$ws.Range("L6").Value = 1.099034652266051
$ws.Range("M6").Value = 1.110548450860428
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.094784852124899
$ws.Range("D7").Value = 1.096917788565025
$ws.Range("E7").Value = 1.096371922944853
$ws.Range("F7").Value = 1.107900986912526
$ws.Range("I7").Value = 1.075749221370079
$ws.Range("J7").Value = 1.099049851313098
$ws.Range("K7").Value = 1.099256220561725
$ws.Range("L7").Value = 1.098711576284521
$ws.Range("M7").Value = 1.11021511798842
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.0930056292286
$ws.Range("D8").Value = 1.095442039989701
$ws.Range("E8").Value = 1.094780759327748
$ws.Range("F8").Value = 1.10627035056697
$ws.Range("I8").Value = 1.075047310197882
$ws.Range("J8").Value = 1.097712414109124
$ws.Range("K8").Value = 1.09801717275362
$ws.Range("L8").Value = 1.097357531201517
$ws.Range("M8").Value = 1.10881894304293
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.089862402254349
$ws.Range("D9").Value = 1.092834277997749
$ws.Range("E9").Value = 1.091967358149237
$ws.Range("F9").Value = 1.103390189093312
$ws.Range("I9").Value = 1.073796739571401
$ws.Range("J9").Value = 1.095343539088143
$ws.Range("K9").Value = 1.095821958852635
$ws.Range("L9").Value = 1.094957557973354
$ws.Range("M9").Value = 1.106347542973437
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.087761235864156
$ws.Range("D10").Value = 1.0910906482932
$ws.Range("E10").Value = 1.090085071870448
$ws.Range("F10").Value = 1.10146528933287
$ws.Range("I10").Value = 1.07295367327229
$ws.Range("J10").Value = 1.09375590215109
$ws.Range("K10").Value = 1.094350310977559
$ws.Range("L10").Value = 1.093347944010548
$ws.Range("M10").Value = 1.104692213637827
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.086849993315685
$ws.Range("D11").Value = 1.090334375030106
$ws.Range("E11").Value = 1.089268378542999
$ws.Range("F11").Value = 1.100630595982352
$ws.Range("I11").Value = 1.072586368100142
$ws.Range("J11").Value = 1.093066397959227
$ws.Range("K11").Value = 1.093711087300124
$ws.Range("L11").Value = 1.092648625495523
$ws.Range("M11").Value = 1.103973553804273
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.086511297544612
$ws.Range("D12").Value = 1.090053266462171
$ws.Range("E12").Value = 1.088964769250034
$ws.Range("F12").Value = 1.100320368829793
$ws.Range("I12").Value = 1.072449593642548
$ws.Range("J12").Value = 1.092809973427754
$ws.Range("K12").Value = 1.093473348106867
$ws.Range("L12").Value = 1.092388510457847
$ws.Range("M12").Value = 1.10370632294054
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.086583958973782
$ws.Range("D13").Value = 1.090113574096143
$ws.Range("E13").Value = 1.089029906004496
$ws.Range("F13").Value = 1.100386922049508
$ws.Range("I13").Value = 1.0724789476998
$ws.Range("J13").Value = 1.092864991549393
$ws.Range("K13").Value = 1.093524357749619
$ws.Range("L13").Value = 1.092444322247125
$ws.Range("M13").Value = 1.103763657999668
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.086822001133555
$ws.Range("D14").Value = 1.090311142510543
$ws.Range("E14").Value = 1.089243287294793
$ws.Range("F14").Value = 1.100604956308066
$ws.Range("I14").Value = 1.072575069251658
$ws.Range("J14").Value = 1.093045208217461
$ws.Range("K14").Value = 1.093691441909858
$ws.Range("L14").Value = 1.092627131615338
$ws.Range("M14").Value = 1.103951470330994
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.086968637475199
$ws.Range("D15").Value = 1.090432844965479
$ws.Range("E15").Value = 1.089374724871778
$ws.Range("F15").Value = 1.100739269770678
$ws.Range("I15").Value = 1.072634247652136
$ws.Range("J15").Value = 1.093156204181448
$ws.Range("K15").Value = 1.09379434769675
$ws.Range("L15").Value = 1.092739719053433
$ws.Range("M15").Value = 1.104067149330471
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.087821681539687
$ws.Range("D16").Value = 1.091140812542667
$ws.Range("E16").Value = 1.09013923790787
$ws.Range("F16").Value = 1.101520659515508
$ws.Range("I16").Value = 1.072978002436918
$ws.Range("J16").Value = 1.093801618768132
$ws.Range("K16").Value = 1.094392691869211
$ws.Range("L16").Value = 1.093394305670274
$ws.Range("M16").Value = 1.104739868496715
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.088356388289696
$ws.Range("D17").Value = 1.091584558910458
$ws.Range("E17").Value = 1.09061835121584
$ws.Range("F17").Value = 1.102010480304643
$ws.Range("I17").Value = 1.073193025974296
$ws.Range("J17").Value = 1.094205919135194
$ws.Range("K17").Value = 1.094767481634906
$ws.Range("L17").Value = 1.093804279423583
$ws.Range("M17").Value = 1.105161337961812
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.088668136726491
$ws.Range("D18").Value = 1.091843266398529
$ws.Range("E18").Value = 1.090897650949365
$ws.Range("F18").Value = 1.102296069114782
$ws.Range("I18").Value = 1.073318228469438
$ws.Range("J18").Value = 1.094441543503309
$ws.Range("K18").Value = 1.094985898389372
$ws.Range("L18").Value = 1.094043184233238
$ws.Range("M18").Value = 1.105406991789387
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.088774411831223
$ws.Range("D19").Value = 1.091931458354767
$ws.Range("E19").Value = 1.090992858129788
$ws.Range("F19").Value = 1.102393428092578
$ws.Range("I19").Value = 1.073360882540902
$ws.Range("J19").Value = 1.094521852009795
$ws.Range("K19").Value = 1.095060340450779
$ws.Range("L19").Value = 1.094124606416492
$ws.Range("M19").Value = 1.105490722629399
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.088299033526399
$ws.Range("D20").Value = 1.09153696180231
$ws.Range("E20").Value = 1.090566963370951
$ws.Range("F20").Value = 1.101957939138369
$ws.Range("I20").Value = 1.073169978458847
$ws.Range("J20").Value = 1.094162561962378
$ws.Range("K20").Value = 1.094727290106422
$ws.Range("L20").Value = 1.093760316525501
$ws.Range("M20").Value = 1.105116137154265
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.086751909762751
$ws.Range("D21").Value = 1.090252968953618
$ws.Range("E21").Value = 1.089180458848225
$ws.Range("F21").Value = 1.10054075579496
$ws.Range("I21").Value = 1.072546773283905
$ws.Range("J21").Value = 1.092992147548363
$ws.Range("K21").Value = 1.093642248187397
$ws.Range("L21").Value = 1.09257330870563
$ws.Range("M21").Value = 1.103896172283289
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.085777897531628
$ws.Range("D22").Value = 1.089444541357288
$ws.Range("E22").Value = 1.088307241991662
$ws.Range("F22").Value = 1.099648644371218
$ws.Range("I22").Value = 1.072152965434763
$ws.Range("J22").Value = 1.092254454344442
$ws.Range("K22").Value = 1.092958283861301
$ws.Range("L22").Value = 1.091824921922994
$ws.Range("M22").Value = 1.103127459954871
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.086294362569605
$ws.Range("D23").Value = 1.089873212668502
$ws.Range("E23").Value = 1.088770291443003
$ws.Range("F23").Value = 1.100121672717035
$ws.Range("I23").Value = 1.072361918425667
$ws.Range("J23").Value = 1.092645692262646
$ws.Range("K23").Value = 1.093321034030583
$ws.Range("L23").Value = 1.092221853495735
$ws.Range("M23").Value = 1.103535128883842
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.088324950097975
$ws.Range("D24").Value = 1.091558469264247
$ws.Range("E24").Value = 1.090590183816959
$ws.Range("F24").Value = 1.101981680586029
$ws.Range("I24").Value = 1.073180393309648
$ws.Range("J24").Value = 1.094182153811533
$ws.Range("K24").Value = 1.09474545152045
$ws.Range("L24").Value = 1.093780182164642
$ws.Range("M24").Value = 1.105136562016212
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.090675981354331
$ws.Range("D25").Value = 1.093509334200745
$ws.Range("E25").Value = 1.092695849543977
$ws.Range("F25").Value = 1.104135606906982
$ws.Range("I25").Value = 1.0741216807572
$ws.Range("J25").Value = 1.09595741176209
$ws.Range("K25").Value = 1.096390899718363
$ws.Range("L25").Value = 1.095579688785506
$ws.Range("M25").Value = 1.10698780490466
